# Apply updated crypto price/volume data as per Mon Jan 15 19:30:18 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.031.08'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.55%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.543.14'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.76%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.87'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.21%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.15'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.56%  '

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.04%  '

# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.10%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.536'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.18%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.43'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.40%  '

# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.16%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.67'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.36%  '

# Row 13
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.02%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.933.84'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.73%  '

# Row 15
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.546.26'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.43%  '

# Row 16
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.46'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.00%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.855'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.43%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.085.11'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.66%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.12'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.93%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.67'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.28%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0974'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.44%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.55'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.01'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.48%  '

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.69%  '

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.29%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.15'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.84%  '

# Row 27
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.02%  '

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.17%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.15'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.84%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.30'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.88%  '

# Row 31
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.97%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '154.68'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.39%  '

# Row 33
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.43%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.19'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.59%  '

# Row 35
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.24%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0795'

# Row 37
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.18%  '

# Row 38
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.97%  '

# Row 39
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.119'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.03%  '

# Row 40
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.01'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.01%  '

# Row 41
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +10.36%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.84'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.92%  '

# Row 43
$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.36'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.56%  '

# Row 44
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.12%  '

# Row 45
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.32%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.028.55'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.18%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.86'

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.83%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '74.97'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.55%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.786.23'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.53%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.74'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.39%  '
